$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new data row (row 90) reported on 2021/10/06, following the
# same layout as the existing time-series rows above it.
#
# Column A stores the report date as literal text (not a real date value)
# formatted with the "yyyy/mm/dd" style already used by the sheet, so we
# build it through a text formula and flatten it back to a static value —
# this avoids Excel's automatic "looks like a date" conversion turning the
# text into a date serial number.
$ws.Range("A90").Formula = "=""2021/10/06"""
$ws.Range("A90").Copy()
$ws.Range("A90").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$excel.CutCopyMode = $false

$ws.Range("B90").Value = 90.5
$ws.Range("C90").Value = 90.9
$ws.Range("D90").Value = 0.91
$ws.Range("E90").Value = 0.9

# Keep the same selection convention used by the source file: the active
# cell tracks one row past the last data row.
$ws.Range("A91").Select()
